$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for 2022-Q1, pushing the rest down
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Column A on the new row needs the same bold/bordered style as the rows below it
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 6.82

# The "index" column (A) is a 0-based row counter; renumber the rows that
# got pushed down by the insert so they stay sequential (1..5).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# ------------------------------------------------------------------
# 2) New "2022-Q1" sheet, positioned right before "总计"
#    Built by copying the "2021-Q4" sheet (same column layout) then
#    overwriting its data with the 2022-Q1 fund-holdings figures.
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($totalSheet)
$newQ = $wb.Worksheets.Item("2021-Q4 (2)")
$newQ.Name = "2022-Q1"

# Expand from 1 data row (row 2) to 6 data rows (rows 2-7)
$newQ.Range("A3:H7").Insert()
$newQ.Range("A2").Copy()
$newQ.Range("A3:A7").PasteSpecial(-4122)
$newQ.Range("B2:H7").ClearFormats()

# Row 2
$newQ.Range("A2").Value = 0
$newQ.Range("B2:G2").NumberFormat = "@"
$newQ.Range("B2").Value = "001071"
$newQ.Range("C2").Value = "华安媒体互联网混合"
$newQ.Range("D2").Value = "51.61"
$newQ.Range("E2").Value = "92.88"
$newQ.Range("F2").Value = "4.89"
$newQ.Range("G2").Value = "2.5237"
$newQ.Range("B2:G2").ClearFormats()
$newQ.Range("H2").Value = 2

# Row 3
$newQ.Range("A3").Value = 1
$newQ.Range("B3:G3").NumberFormat = "@"
$newQ.Range("B3").Value = "001694"
$newQ.Range("C3").Value = "华安沪港深外延增长灵活配置混合"
$newQ.Range("D3").Value = "43.58"
$newQ.Range("E3").Value = "92.63"
$newQ.Range("F3").Value = "4.80"
$newQ.Range("G3").Value = "2.0918"
$newQ.Range("B3:G3").ClearFormats()
$newQ.Range("H3").Value = 4

# Row 4
$newQ.Range("A4").Value = 2
$newQ.Range("B4:G4").NumberFormat = "@"
$newQ.Range("B4").Value = "006879"
$newQ.Range("C4").Value = "华安智能生活混合"
$newQ.Range("D4").Value = "33.48"
$newQ.Range("E4").Value = "92.78"
$newQ.Range("F4").Value = "4.71"
$newQ.Range("G4").Value = "1.5769"
$newQ.Range("B4:G4").ClearFormats()
$newQ.Range("H4").Value = 4

# Row 5
$newQ.Range("A5").Value = 3
$newQ.Range("B5:G5").NumberFormat = "@"
$newQ.Range("B5").Value = "007460"
$newQ.Range("C5").Value = "华安成长创新混合"
$newQ.Range("D5").Value = "13.21"
$newQ.Range("E5").Value = "91.10"
$newQ.Range("F5").Value = "4.70"
$newQ.Range("G5").Value = "0.6209"
$newQ.Range("B5:G5").ClearFormats()
$newQ.Range("H5").Value = 5

# Row 6
$newQ.Range("A6").Value = 4
$newQ.Range("B6:G6").NumberFormat = "@"
$newQ.Range("B6").Value = "006346"
$newQ.Range("C6").Value = "安信量化优选股票A"
$newQ.Range("D6").Value = "0.71"
$newQ.Range("E6").Value = "90.62"
$newQ.Range("F6").Value = "0.60"
$newQ.Range("G6").Value = "0.0043"
$newQ.Range("B6:G6").ClearFormats()
$newQ.Range("H6").Value = 10

# Row 7
$newQ.Range("A7").Value = 5
$newQ.Range("B7:G7").NumberFormat = "@"
$newQ.Range("B7").Value = "006347"
$newQ.Range("C7").Value = "安信量化优选股票C"
$newQ.Range("D7").Value = "0.49"
$newQ.Range("E7").Value = "90.62"
$newQ.Range("F7").Value = "0.60"
$newQ.Range("G7").Value = "0.0029"
$newQ.Range("B7:G7").ClearFormats()
$newQ.Range("H7").Value = 10

